$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Updated symbol list (price + 1h volume change columns) for rows 2-27 and 39-47.
# Leading apostrophe forces these numeric-looking strings to be stored as text
# (matching the original inlineStr text cells) instead of being parsed as
# numbers/percentages by Excel.
$updates = [ordered]@{
    "D2"  = "'301.01";  "E2"  = "'-0.78%"
    "D3"  = "'31.39";   "E3"  = "'-1.71%"
    "D4"  = "'5.145";   "E4"  = "'-2.29%"
    "D5"  = "'0.07382"; "E5"  = "'-1.26%"
    "D6"  = "'2.519";   "E6"  = "'65.65%"
    "D7"  = "'7.941";   "E7"  = "'1.33%"
    "D8"  = "'3.765";   "E8"  = "'-0.85%"
    "D9"  = "'0.9206";  "E9"  = "'0.28%"
    "D10" = "'0.1735";  "E10" = "'3.15%"
    "D11" = "'0.07556"; "E11" = "'-5.28%"
    "D12" = "'0.08132"; "E12" = "'1.42%"
    "D13" = "'0.03035"; "E13" = "'1.20%"
    "D14" = "'0.09916"; "E14" = "'0.25%"
    "D15" = "'0.001494";"E15" = "'-0.23%"
    "D16" = "'0.006104";"E16" = "'-5.76%"
    "D17" = "'3.459";   "E17" = "'-0.22%"
    "E18" = "'-0.20%"
    "D19" = "'0.3279";  "E19" = "'-1.41%"
    "D20" = "'0.1337";  "E20" = "'0.01%"
    "D21" = "'4.647";   "E21" = "'3.61%"
    "D22" = "'0.04639"; "E22" = "'0.88%"
    "E23" = "'-3.25%"
    "D24" = "'0.001222";"E24" = "'0.38%"
    "D25" = "'0.004493";"E25" = "'1.09%"
    "D26" = "'0.0001299";"E26" = "'-7.12%"
    "E27" = "'5.38%"
    "D39" = "'0.01728"; "E39" = "'0.50%"
    "D40" = "'0.04522"; "E40" = "'0.63%"
    "D41" = "'0.007180";"E41" = "'0.38%"
    "D42" = "'0.1345";  "E42" = "'-0.20%"
    "D43" = "'0.002227";"E43" = "'3.26%"
    "D44" = "'0.01075"; "E44" = "'-16.19%"
    "D45" = "'0.00006275";"E45" = "'1.51%"
    "D46" = "'1.928";   "E46" = "'3.21%"
    "E47" = "'-23.01%"
}

foreach ($addr in $updates.Keys) {
    $ws.Range($addr).Value = $updates[$addr]
}
